$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete rows 6 and 7 (Resolving-Mac sending cluster rows no longer present)
$ws.Rows("6:7").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Gfra3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.534048
$ws.Range("H2").Value = 1.602144
$ws.Range("I2").Value = 0.2492808729834395
$ws.Range("J2").Value = 0.3324807621550537
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.036591
$ws.Range("N2").Value = 0.073182
$ws.Range("O2").Value = 0.01193099839087543
$ws.Range("P2").Value = 0.01193099839087543
$ws.Range("Q2").Value = 0.019541350368
$ws.Range("R2").Value = 0.117248102208
$ws.Range("S2").Value = 0.002974169694441439
$ws.Range("T2").Value = 0.003966827438268982

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Gfra3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.534048
$ws.Range("H3").Value = 1.602144
$ws.Range("I3").Value = 0.2492808729834395
$ws.Range("J3").Value = 0.3324807621550537
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.030294
$ws.Range("N3").Value = 6.060588
$ws.Range("O3").Value = 0.9880690016091246
$ws.Range("P3").Value = 0.9880690016091246
$ws.Range("Q3").Value = 1.618322450112
$ws.Range("R3").Value = 9.709934700672001
$ws.Range("S3").Value = 0.246306703288998
$ws.Range("T3").Value = 0.3285139347167847

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Gfra3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.6083065
$ws.Range("H4").Value = 3.216613
$ws.Range("I4").Value = 0.7507191270165605
$ws.Range("J4").Value = 0.6675192378449464
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.036591
$ws.Range("N4").Value = 0.073182
$ws.Range("O4").Value = 0.01193099839087543
$ws.Range("P4").Value = 0.01193099839087543
$ws.Range("Q4").Value = 0.0588495431415
$ws.Range("R4").Value = 0.235398172566
$ws.Range("S4").Value = 0.008956828696433992
$ws.Range("T4").Value = 0.00796417095260645

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Gfra3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.6083065
$ws.Range("H5").Value = 3.216613
$ws.Range("I5").Value = 0.7507191270165605
$ws.Range("J5").Value = 0.6675192378449464
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.030294
$ws.Range("N5").Value = 6.060588
$ws.Range("O5").Value = 0.9880690016091246
$ws.Range("P5").Value = 0.9880690016091246
$ws.Range("Q5").Value = 4.873641537111
$ws.Range("R5").Value = 19.494566148444
$ws.Range("S5").Value = 0.7417622983201265
$ws.Range("T5").Value = 0.6595550668923399
